# "updated 2 icdc scripts to resolve wait time issue"
#
# The "CasesTab" Cypher query stored in B2 of the "startup" sheet dropped its
# trailing OPTIONAL MATCH (co:cohort) clause and the `Cohort` output column
# (removed to cut down the query's run time). Apply that same text edit here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("B2")
$oldValue = $cell.Value2

$cohortSuffixCRLF = ",`r`n        coalesce(co.cohort_description, '') AS ``Cohort``"
$cohortSuffixLF   = ",`n        coalesce(co.cohort_description, '') AS ``Cohort``"

if ($oldValue.Contains($cohortSuffixCRLF)) {
    $newValue = $oldValue.Replace($cohortSuffixCRLF, "")
} else {
    $newValue = $oldValue.Replace($cohortSuffixLF, "")
}

$cell.Value2 = $newValue

# The query text lost a line, so the wrapped row shrinks accordingly.
$ws.Rows.Item(2).RowHeight = 304.5

# Reflect the edited cell as the new selection / scroll position.
$cell.Select()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
